$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-12) get rearranged: the full record (D, L, M, N, O, P, R, S)
# currently sitting in one row moves to another row, per the mapping below
# (new row -> old row it takes its values from).
$mapping = @{
    2  = 5
    3  = 6
    4  = 7
    5  = 8
    6  = 11
    7  = 12
    8  = 9
    9  = 10
    10 = 2
    11 = 3
    12 = 4
}

$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

# Snapshot the current values for every relevant cell before writing anything,
# since several rows swap values with each other.
$snapshot = @{}
foreach ($row in 2..12) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$oldRow"
        $dstAddr = "$col$newRow"
        $ws.Range($dstAddr).Value2 = $snapshot[$srcAddr]
    }
}
